# Converted Fig. 4->5, Fig.S1->4
# The data blocks that used to be labelled "Figure 4 [...]" (rows 64-86)
# and "Figure S1... " (rows 89-111) are swapped in place (same size,
# 23 rows x 8 cols each), and the section titles are relabelled:
#   old Figure 4 [...]  -> new Figure 5 [...]   (content moves to 89-111)
#   old Figure S1...     -> new Figure 4...      (content moves to 64-86)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcFig4  = $ws.Range("A64:H86")
$srcFigS1 = $ws.Range("A89:H111")
$scratch  = $ws.Range("A500:H522")

# 1) stash the "Figure 4" block out of the way
$srcFig4.Copy($scratch)

# 2) move "Figure S1" block into the old "Figure 4" rows
$srcFigS1.Copy($srcFig4)

# 3) move the stashed "Figure 4" block into the old "Figure S1" rows
$scratch.Copy($srcFigS1)

# 4) clear the scratch area
$scratch.Clear()

# 5) relabel the section titles to their new figure numbers
$ws.Range("A64").Value2  = "Figure 4B (E, F, G, H)"
$ws.Range("A70").Value2  = "Figure 4C"
$ws.Range("A76").Value2  = "Figure 4D"
$ws.Range("A82").Value2  = "Figure 4 (I & J)"

$ws.Range("A89").Value2  = "Figure 5 [Correlation: Drug Condition]"
$ws.Range("A95").Value2  = "Figure 5 [Correlation: FR Control]"
$ws.Range("A101").Value2 = "Figure 5 [MI Value: Drug Condition]"
$ws.Range("A107").Value2 = "Figure 5 [MI Value: FR Control]"

# 6) keep the workbook's Excel Tables (ListObjects) anchored to the data
#    that now lives in their (possibly new) rows
$ws.ListObjects.Item("Table10").Resize($ws.Range("A90:D93"))
$ws.ListObjects.Item("Table12").Resize($ws.Range("A96:B99"))
$ws.ListObjects.Item("Table13").Resize($ws.Range("A102:D105"))
$ws.ListObjects.Item("Table14").Resize($ws.Range("A108:B111"))
$ws.ListObjects.Item("Table15").Resize($ws.Range("A65:D68"))
$ws.ListObjects.Item("Table16").Resize($ws.Range("A71:D74"))
$ws.ListObjects.Item("Table17").Resize($ws.Range("A77:H80"))
$ws.ListObjects.Item("Table18").Resize($ws.Range("A83:B86"))

# 7) reflect the scrolled/selected view seen in the edited workbook
$ws.Application.ActiveWindow.ScrollRow = 84
$ws.Range("A107").Select()
